{"js": "// 1) \"Data de entrega: 08/11\" -> bold the whole line and change the date\n//    to \"09/05/2025\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load('items/text');\nawait context.sync();\n\nlet dateParaIndex = -1;\nlet obsParaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (dateParaIndex === -1 && t.indexOf('Data de entrega') !== -1) {\n    dateParaIndex = i;\n  }\n  if (obsParaIndex === -1 && t.indexOf('Como s\u00e3o 32 alunos') !== -1) {\n    obsParaIndex = i;\n  }\n}\nif (dateParaIndex === -1) {\n  throw new Error('Could not find the \"Data de entrega\" paragraph.');\n}\nif (obsParaIndex === -1) {\n  throw new Error('Could not find the \"OBS.\" paragraph.');\n}\n\nconst dateParagraph = paragraphs.items[dateParaIndex];\n// Bold the paragraph mark too, matching the rest of the (already bold) line.\ndateParagraph.font.bold = true;\n\nconst dateResults = dateParagraph.search('08/11', { matchCase: true });\ndateResults.load('items');\nawait context.sync();\nif (dateResults.items.length === 0) {\n  throw new Error('Could not find \"08/11\" text to replace.');\n}\ndateResults.items[0].insertText('09/05/2025', Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Remove the sentence about group sizes / counts from the OBS paragraph.\nconst obsParagraph = paragraphs.items[obsParaIndex];\nconst removedText =\n  'Como s\u00e3o 32 alunos no total, ser\u00e3o 8 grupos de 4 pessoas. ' +\n  'Os grupos poder\u00e3o ser compostos por integrantes de turmas diferentes. ';\nconst obsResults = obsParagraph.search(removedText, { matchCase: true });\nobsResults.load('items');\nawait context.sync();\nif (obsResults.items.length === 0) {\n  throw new Error('Could not find the sentence to remove.');\n}\nobsResults.items[0].insertText('', Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"Data de entrega: 08/11\" -> bold the whole line and change the date\n#    to \"09/05/2025\".\n$dateParagraph = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*Data de entrega*\") {\n        $dateParagraph = $d.Paragraphs.Item($i)\n        break\n    }\n}\nif ($dateParagraph -eq $null) {\n    throw \"Could not find the 'Data de entrega' paragraph.\"\n}\n\n# Bold the paragraph mark too, matching the rest of the (already bold) line.\n$dateParagraph.Range.Font.Bold = 1\n\n$dateRange = $dateParagraph.Range\n# 1 = wdFindContinue, 2 = wdReplaceOne\n$found = $dateRange.Find.Execute(\"08/11\", $false, $false, $false, $false, $false, $true, 1, $false, \"09/05/2025\", 2)\nif (-not $found) {\n    throw \"Could not find '08/11' text to replace.\"\n}\n\n# 2) Remove the sentence about group sizes / counts from the OBS paragraph.\n$obsParagraph = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*Como s\u00e3o 32 alunos*\") {\n        $obsParagraph = $d.Paragraphs.Item($i)\n        break\n    }\n}\nif ($obsParagraph -eq $null) {\n    throw \"Could not find the 'OBS.' paragraph.\"\n}\n\n$removedText = \"Como s\u00e3o 32 alunos no total, ser\u00e3o 8 grupos de 4 pessoas. Os grupos poder\u00e3o ser compostos por integrantes de turmas diferentes. \"\n$obsRange = $obsParagraph.Range\n$found2 = $obsRange.Find.Execute($removedText, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\nif (-not $found2) {\n    throw \"Could not find the sentence to remove.\"\n}\n"}
